$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = 'D2'; Value = '30.130.72' },
    @{ Cell = 'E2'; Value = '  -1.25%  ' },
    @{ Cell = 'D3'; Value = '1.867.76' },
    @{ Cell = 'E3'; Value = '  +0.86%  ' },
    @{ Cell = 'D4'; Value = '0.9935' },
    @{ Cell = 'E4'; Value = '  -0.72%  ' },
    @{ Cell = 'D5'; Value = '236.33' },
    @{ Cell = 'E5'; Value = '  +1.10%  ' },
    @{ Cell = 'D6'; Value = '0.9942' },
    @{ Cell = 'E6'; Value = '  -0.64%  ' },
    @{ Cell = 'D7'; Value = '0.4655' },
    @{ Cell = 'E7'; Value = '  -1.19%  ' },
    @{ Cell = 'D8'; Value = '0.2818' },
    @{ Cell = 'E8'; Value = '  +2.80%  ' },
    @{ Cell = 'D9'; Value = '0.06475' },
    @{ Cell = 'E9'; Value = '  +2.13%  ' },
    @{ Cell = 'B10'; Value = 'Litecoin' },
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' },
    @{ Cell = 'D10'; Value = '107.26' },
    @{ Cell = 'E10'; Value = '  +26.80%  ' },
    @{ Cell = 'B11'; Value = 'Solana' },
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol' },
    @{ Cell = 'D11'; Value = '18.74' },
    @{ Cell = 'E11'; Value = '  +6.23%  ' },
    @{ Cell = 'D12'; Value = '1.844.11' },
    @{ Cell = 'E12'; Value = '  -0.37%  ' },
    @{ Cell = 'D13'; Value = '0.07520' },
    @{ Cell = 'E13'; Value = '  +1.51%  ' },
    @{ Cell = 'D14'; Value = '5.030' },
    @{ Cell = 'E14'; Value = '  -0.46%  ' },
    @{ Cell = 'B15'; Value = 'BitcoinCash' },
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' },
    @{ Cell = 'D15'; Value = '312.55' },
    @{ Cell = 'E15'; Value = '  +29.13%  ' },
    @{ Cell = 'B16'; Value = 'Polygon' },
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' },
    @{ Cell = 'D16'; Value = '0.6341' },
    @{ Cell = 'E16'; Value = '  +1.40%  ' },
    @{ Cell = 'D17'; Value = '30.100.70' },
    @{ Cell = 'E17'; Value = '  -1.22%  ' },
    @{ Cell = 'D18'; Value = '0.9979' },
    @{ Cell = 'E18'; Value = '  -0.24%  ' },
    @{ Cell = 'D19'; Value = '12.71' },
    @{ Cell = 'E19'; Value = '  +0.29%  ' },
    @{ Cell = 'D20'; Value = '0.000007440' },
    @{ Cell = 'E20'; Value = '  +1.37%  ' },
    @{ Cell = 'D21'; Value = '2.083.50' },
    @{ Cell = 'E21'; Value = '  -0.20%  ' },
    @{ Cell = 'D22'; Value = '0.9909' },
    @{ Cell = 'E22'; Value = '  -1.01%  ' },
    @{ Cell = 'D23'; Value = '5.041' },
    @{ Cell = 'E23'; Value = '  +2.13%  ' },
    @{ Cell = 'D24'; Value = '6.175' },
    @{ Cell = 'E24'; Value = '  +3.40%  ' },
    @{ Cell = 'B25'; Value = 'Cosmos' },
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' },
    @{ Cell = 'D25'; Value = '9.151' },
    @{ Cell = 'E25'; Value = '  -0.86%  ' },
    @{ Cell = 'B26'; Value = 'Monero' },
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' },
    @{ Cell = 'D26'; Value = '164.76' },
    @{ Cell = 'E26'; Value = '  +1.79%  ' },
    @{ Cell = 'D27'; Value = '19.91' },
    @{ Cell = 'E27'; Value = '  +10.56%  ' },
    @{ Cell = 'D28'; Value = '1.967' },
    @{ Cell = 'E28'; Value = '  +4.44%  ' },
    @{ Cell = 'E29'; Value = '  +6.20%  ' },
    @{ Cell = 'E30'; Value = '  -2.82%  ' },
    @{ Cell = 'D31'; Value = '4.041' },
    @{ Cell = 'E31'; Value = '  +0.58%  ' },
    @{ Cell = 'D32'; Value = '3.879' },
    @{ Cell = 'E32'; Value = '  +0.93%  ' },
    @{ Cell = 'D33'; Value = '0.04906' },
    @{ Cell = 'E33'; Value = '  +0.81%  ' },
    @{ Cell = 'D34'; Value = '0.7445' },
    @{ Cell = 'E34'; Value = '  +5.52%  ' },
    @{ Cell = 'D35'; Value = '1.134' },
    @{ Cell = 'E35'; Value = '  -0.21%  ' },
    @{ Cell = 'D36'; Value = '2.708' },
    @{ Cell = 'E36'; Value = '  -0.17%  ' },
    @{ Cell = 'D37'; Value = '0.01922' },
    @{ Cell = 'E37'; Value = '  +1.26%  ' },
    @{ Cell = 'D38'; Value = '2.650' },
    @{ Cell = 'E38'; Value = '  -1.41%  ' },
    @{ Cell = 'D39'; Value = '1.991' },
    @{ Cell = 'E39'; Value = '  +0.82%  ' },
    @{ Cell = 'B40'; Value = 'TrustWalletToken' },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' },
    @{ Cell = 'D40'; Value = '0.8631' },
    @{ Cell = 'E40'; Value = '  -1.22%  ' },
    @{ Cell = 'B41'; Value = 'Quant' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' },
    @{ Cell = 'D41'; Value = '107.06' },
    @{ Cell = 'E41'; Value = '  +1.65%  ' },
    @{ Cell = 'B42'; Value = 'FraxShare' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' },
    @{ Cell = 'D42'; Value = '5.773' },
    @{ Cell = 'E42'; Value = '  +4.88%  ' },
    @{ Cell = 'B43'; Value = 'PaxDollar' },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' },
    @{ Cell = 'D43'; Value = '0.9939' },
    @{ Cell = 'E43'; Value = '  -0.66%  ' },
    @{ Cell = 'D44'; Value = '0.4088' },
    @{ Cell = 'E44'; Value = '  +0.35%  ' },
    @{ Cell = 'D45'; Value = '66.85' },
    @{ Cell = 'E45'; Value = '  +7.34%  ' },
    @{ Cell = 'D46'; Value = '7.154' },
    @{ Cell = 'E46'; Value = '  -0.77%  ' },
    @{ Cell = 'D47'; Value = '9.161' },
    @{ Cell = 'E47'; Value = '  +7.17%  ' },
    @{ Cell = 'D48'; Value = '0.1197' },
    @{ Cell = 'D49'; Value = '34.13' },
    @{ Cell = 'E49'; Value = '  +2.43%  ' },
    @{ Cell = 'D50'; Value = '0.05553' },
    @{ Cell = 'E50'; Value = '  +0.28%  ' },
    @{ Cell = 'D51'; Value = '0.3760' },
    @{ Cell = 'E51'; Value = '  +2.22%  ' }
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $change.Value
    $cell.ClearFormats()
}

